$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both contain the same event listing and both
# need their "想去人数" (want-to-go count) figures bumped.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 191
    $ws.Range("F4").Value = 142
}
